$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -85.714285714285
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 22
$ws.Range("K15").Value = -50
$ws.Range("N15").Value = -86.25
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -29.411764705882
$ws.Range("I16").Value = 98
$ws.Range("J16").Value = 127
$ws.Range("K16").Value = -22.834645669291
$ws.Range("L16").Value = 4.255319148936
$ws.Range("M16").Value = -57.017543859649
$ws.Range("N16").Value = -92.530487804878
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 24.137931034482
$ws.Range("I17").Value = 208
$ws.Range("J17").Value = 230
$ws.Range("K17").Value = -9.565217391304
$ws.Range("L17").Value = 9.473684210526
$ws.Range("M17").Value = -9.565217391304
$ws.Range("N17").Value = -62.857142857142
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -63.636363636363
$ws.Range("I18").Value = 76
$ws.Range("J18").Value = 113
$ws.Range("K18").Value = -32.743362831858
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -55.555555555555
$ws.Range("N18").Value = -95.583962812318
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -46.153846153846
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -18.60465116279
$ws.Range("I19").Value = 325
$ws.Range("J19").Value = 322
$ws.Range("K19").Value = 0.931677018633
$ws.Range("L19").Value = 13.636363636363
$ws.Range("M19").Value = -16.879795396419
$ws.Range("N19").Value = -50.906344410876
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 65
$ws.Range("J20").Value = 79
$ws.Range("K20").Value = -17.721518987341
$ws.Range("L20").Value = -18.75
$ws.Range("M20").Value = -46.280991735537
$ws.Range("N20").Value = -95.849297573435
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -6.896551724137
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 125
$ws.Range("H21").Value = -19.2
$ws.Range("I21").Value = 787
$ws.Range("J21").Value = 894
$ws.Range("K21").Value = -11.968680089485
$ws.Range("L21").Value = 0.897435897435
$ws.Range("M21").Value = -32.330180567497
$ws.Range("N21").Value = -86.726260752234
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -50
$ws.Range("M22").Value = -21.052631578947
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = -12.195121951219
$ws.Range("F24").Value = 143
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = 22.222222222222
$ws.Range("I24").Value = 1145
$ws.Range("J24").Value = 819
$ws.Range("K24").Value = 39.804639804639
$ws.Range("L24").Value = 42.590286425902
$ws.Range("M24").Value = 45.120405576679
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 118.181818181818
$ws.Range("F25").Value = 56
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = 27.272727272727
$ws.Range("I25").Value = 378
$ws.Range("J25").Value = 368
$ws.Range("K25").Value = 2.717391304347
$ws.Range("L25").Value = 17.027863777089
$ws.Range("M25").Value = -21.085594989561
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("J26").Value = 31
$ws.Range("K26").Value = -35.483870967741
$ws.Range("L26").Value = -23.076923076923
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -54.545454545454
$ws.Range("I27").Value = 42
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = -8.695652173913
$ws.Range("L27").Value = 13.513513513513
$ws.Range("N28").Value = -90.322580645161
$ws.Range("N29").Value = -89.189189189189
$ws.Range("I30").Value = 11
$ws.Range("K30").Value = 10
$ws.Range("L30").Value = 120

# --- Type-changing cells: text placeholder <-> numeric value ---
# Row 15: D15,E15 go from text ("0","***.*") to numbers (1,-100)
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -100

# Row 26: D26,E26 go from text to numbers (1,-100)
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E26").Value = -100

# Row 30: C30 goes from text to number (1)
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("C30").Value = 1

# Row 22: D22,E22 go from numbers to text ("0","***.*")
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").PasteSpecial(-4122)

# Row 28: D28,E28 go from numbers to text
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)

# Row 29: D29,E29 go from numbers to text
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("C29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)
